$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (46061 -> 46062) for every data row (rows 2 through 437).
$lastRow = 437
$range = $ws.Range("C2:C$lastRow")
$range.Value = 46062
